$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Cells.Item(4, 2).Value = 6469236
$ws.Cells.Item(4, 3).Value = 8986
$ws.Cells.Item(4, 4).Value = 3727743
$ws.Cells.Item(4, 5).Value = 2548188
$ws.Cells.Item(4, 7).Value = 57
$ws.Cells.Item(4, 8).Value = 193305

# Row 5
$ws.Cells.Item(5, 2).Value = 4251587
$ws.Cells.Item(5, 3).Value = 49025
$ws.Cells.Item(5, 4).Value = 3298881
$ws.Cells.Item(5, 5).Value = 880385
$ws.Cells.Item(5, 7).Value = 634
$ws.Cells.Item(5, 8).Value = 72321

# Row 22
$ws.Cells.Item(22, 2).Value = 278784
$ws.Cells.Item(22, 3).Value = 1108
$ws.Cells.Item(22, 4).Value = 210238
$ws.Cells.Item(22, 5).Value = 32993
$ws.Cells.Item(22, 7).Value = 12
$ws.Cells.Item(22, 8).Value = 35553

# Row 24
$ws.Cells.Item(24, 2).Value = 252964
$ws.Cells.Item(24, 3).Value = 1240
$ws.Cells.Item(24, 5).Value = 16563

# Row 28
$ws.Cells.Item(28, 1).Value = "Canada"
$ws.Cells.Item(28, 2).Value = 132111
$ws.Cells.Item(28, 3).Value = 216
$ws.Cells.Item(28, 4).Value = 116446
$ws.Cells.Item(28, 5).Value = 6519
$ws.Cells.Item(28, 7).Value = 1
$ws.Cells.Item(28, 8).Value = 9146

# Row 29
$ws.Cells.Item(29, 1).Value = "Israel"
$ws.Cells.Item(29, 2).Value = 131970
$ws.Cells.Item(29, 3).Value = 1326
$ws.Cells.Item(29, 4).Value = 103849
$ws.Cells.Item(29, 5).Value = 27099
$ws.Cells.Item(29, 7).Value = 3
$ws.Cells.Item(29, 8).Value = 1022

# Row 59
$ws.Cells.Item(59, 2).Value = 46653
$ws.Cells.Item(59, 3).Value = 289
$ws.Cells.Item(59, 4).Value = 32985
$ws.Cells.Item(59, 5).Value = 12106
$ws.Cells.Item(59, 7).Value = 6
$ws.Cells.Item(59, 8).Value = 1562

# Row 60
$ws.Cells.Item(60, 1).Value = "Ghana"
$ws.Cells.Item(60, 2).Value = 44869
$ws.Cells.Item(60, 3).Value = 92
$ws.Cells.Item(60, 4).Value = 43801
$ws.Cells.Item(60, 5).Value = 785
$ws.Cells.Item(60, 7).Value = 0
$ws.Cells.Item(60, 8).Value = 283

# Row 61
$ws.Cells.Item(61, 1).Value = "Armenia"
$ws.Cells.Item(61, 2).Value = 44845
$ws.Cells.Item(61, 3).Value = 62
$ws.Cells.Item(61, 4).Value = 40121
$ws.Cells.Item(61, 5).Value = 3824
$ws.Cells.Item(61, 7).Value = 3
$ws.Cells.Item(61, 8).Value = 900

# Row 72
$ws.Cells.Item(72, 2).Value = 28374
$ws.Cells.Item(72, 3).Value = 218
$ws.Cells.Item(72, 5).Value = 8885

# Row 91
$ws.Cells.Item(91, 5).Value = 7431
$ws.Cells.Item(91, 7).Value = 5
$ws.Cells.Item(91, 8).Value = 289

# Row 139
$ws.Cells.Item(139, 1).Value = "Jordania"
$ws.Cells.Item(139, 2).Value = 2478
$ws.Cells.Item(139, 3).Value = 67
$ws.Cells.Item(139, 4).Value = 1817
$ws.Cells.Item(139, 5).Value = 644
$ws.Cells.Item(139, 7).Value = 1
$ws.Cells.Item(139, 8).Value = 17

# Row 140
$ws.Cells.Item(140, 1).Value = "Aruba"
$ws.Cells.Item(140, 2).Value = 2449
$ws.Cells.Item(140, 4).Value = 1206
$ws.Cells.Item(140, 5).Value = 1229
$ws.Cells.Item(140, 8).Value = 14

# Row 153
$ws.Cells.Item(153, 1).Value = "Birmania"
$ws.Cells.Item(153, 2).Value = 1518
$ws.Cells.Item(153, 3).Value = 99
$ws.Cells.Item(153, 4).Value = 388
$ws.Cells.Item(153, 5).Value = 1122
$ws.Cells.Item(153, 8).Value = 8

# Row 154
$ws.Cells.Item(154, 1).Value = "Republica de Chipre"
$ws.Cells.Item(154, 2).Value = 1509
$ws.Cells.Item(154, 4).Value = 1237
$ws.Cells.Item(154, 5).Value = 251
$ws.Cells.Item(154, 8).Value = 21

# Row 155
$ws.Cells.Item(155, 1).Value = "Togo"
$ws.Cells.Item(155, 2).Value = 1488
$ws.Cells.Item(155, 4).Value = 1106
$ws.Cells.Item(155, 5).Value = 350
$ws.Cells.Item(155, 8).Value = 32

# Row 156
$ws.Cells.Item(156, 1).Value = "Guyana"
$ws.Cells.Item(156, 2).Value = 1468
$ws.Cells.Item(156, 4).Value = 954
$ws.Cells.Item(156, 5).Value = 468
$ws.Cells.Item(156, 8).Value = 46

# Row 157
$ws.Cells.Item(157, 1).Value = "Burkina Faso"
$ws.Cells.Item(157, 2).Value = 1452
$ws.Cells.Item(157, 3).Value = 0
$ws.Cells.Item(157, 4).Value = 1103
$ws.Cells.Item(157, 5).Value = 294
$ws.Cells.Item(157, 8).Value = 55

# Row 158
$ws.Cells.Item(158, 1).Value = "Letonia"
$ws.Cells.Item(158, 2).Value = 1429
$ws.Cells.Item(158, 3).Value = 1
$ws.Cells.Item(158, 4).Value = 1187
$ws.Cells.Item(158, 5).Value = 207
$ws.Cells.Item(158, 8).Value = 35
